# Screenshot capture method has been moved to Generic Functions
#
# Update the Order# test value on the STS_TestData sheet (cell H2) from
# "88016008" to "88016078".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STS_TestData")

$cell = $ws.Range("H2")

# Leading apostrophe forces text entry so the numeric-looking order
# number isn't coerced into a Number type.
$cell.Value = "'88016078"

# The apostrophe trick nudges the cell onto a distinct (text) style
# record, so restore the original formatting by copying it back from a
# neighboring cell on the same row that already carries the correct
# style.
$ws.Range("G2").Copy()
$cell.PasteSpecial(-4122)
$excel.CutCopyMode = $false
